# This script updates the weekly price data (Fruta / hortaliza, semanal) for
# "Hortaliza, Mapocho Venta Directa de Santiago - Ají", rows 2-18.
# Each data row's values (date, variety, quality, volume, prices, unit,
# price per kg, and kg/units) are rotated to a different week's figures,
# per the upstream weekly refresh of this subset workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    # Row 2 (was data row now matching original row 6)
    $ws.Range("D2").Value = 44581
    $ws.Range("I2").Value = 'Segunda'
    $ws.Range("J2").Value = 30
    $ws.Range("K2").Value = 17000
    $ws.Range("L2").Value = 17000
    $ws.Range("M2").Value = 17000
    $ws.Range("P2").Value = 680

    # Row 3 (was data row now matching original row 7)
    $ws.Range("D3").Value = 44449
    $ws.Range("H3").Value = 'Americana (o)'
    $ws.Range("J3").Value = 25
    $ws.Range("N3").Value = '$/caja 25 kilos'
    $ws.Range("P3").Value = 3200
    $ws.Range("Q3").Value = 25

    # Row 4 (was data row now matching original row 8)
    $ws.Range("D4").Value = 44449
    $ws.Range("I4").Value = 'Segunda'
    $ws.Range("J4").Value = 20
    $ws.Range("N4").Value = '$/caja 15 kilos'
    $ws.Range("P4").Value = 5000
    $ws.Range("Q4").Value = 15

    # Row 5 (was data row now matching original row 13)
    $ws.Range("D5").Value = 44193
    $ws.Range("J5").Value = 15
    $ws.Range("K5").Value = 46000
    $ws.Range("L5").Value = 46000
    $ws.Range("M5").Value = 46000
    $ws.Range("N5").Value = '$/caja 15 kilos'
    $ws.Range("P5").Value = 3067
    $ws.Range("Q5").Value = 15

    # Row 6 (was data row now matching original row 17)
    $ws.Range("D6").Value = 44326
    $ws.Range("I6").Value = 'Primera'
    $ws.Range("J6").Value = 15
    $ws.Range("K6").Value = 30000
    $ws.Range("L6").Value = 30000
    $ws.Range("M6").Value = 30000
    $ws.Range("P6").Value = 1200

    # Row 7 (was data row now matching original row 2)
    $ws.Range("D7").Value = 44446
    $ws.Range("J7").Value = 5
    $ws.Range("K7").Value = 78000
    $ws.Range("L7").Value = 78000
    $ws.Range("M7").Value = 78000
    $ws.Range("P7").Value = 3120

    # Row 8 (was data row now matching original row 3)
    $ws.Range("D8").Value = 44446
    $ws.Range("H8").Value = 'Inferno'
    $ws.Range("I8").Value = 'Primera'
    $ws.Range("J8").Value = 4
    $ws.Range("K8").Value = 80000
    $ws.Range("L8").Value = 80000
    $ws.Range("M8").Value = 80000
    $ws.Range("P8").Value = 5333

    # Row 9 (was data row now matching original row 5)
    $ws.Range("D9").Value = 44474
    $ws.Range("J9").Value = 18
    $ws.Range("K9").Value = 100000
    $ws.Range("L9").Value = 100000
    $ws.Range("M9").Value = 100000
    $ws.Range("P9").Value = 4000

    # Row 10 (was data row now matching original row 9)
    $ws.Range("D10").Value = 44319
    $ws.Range("H10").Value = 'Americana (o)'
    $ws.Range("J10").Value = 20
    $ws.Range("K10").Value = 30000
    $ws.Range("L10").Value = 30000
    $ws.Range("M10").Value = 30000
    $ws.Range("P10").Value = 1200

    # Row 11 (was data row now matching original row 15)
    $ws.Range("D11").Value = 44221
    $ws.Range("H11").Value = 'Americana (o)'
    $ws.Range("J11").Value = 22
    $ws.Range("K11").Value = 24000
    $ws.Range("L11").Value = 25000
    $ws.Range("M11").Value = 24545
    $ws.Range("P11").Value = 982

    # Row 12 (was data row now matching original row 11)
    $ws.Range("D12").Value = 44553
    $ws.Range("H12").Value = 'Inferno'
    $ws.Range("J12").Value = 35
    $ws.Range("K12").Value = 45000
    $ws.Range("L12").Value = 45000
    $ws.Range("M12").Value = 45000
    $ws.Range("P12").Value = 1800

    # Row 13 (was data row now matching original row 18)
    $ws.Range("D13").Value = 44343
    $ws.Range("J13").Value = 20
    $ws.Range("K13").Value = 36000
    $ws.Range("L13").Value = 36000
    $ws.Range("M13").Value = 36000
    $ws.Range("N13").Value = '$/caja 25 kilos'
    $ws.Range("P13").Value = 1440
    $ws.Range("Q13").Value = 25

    # Row 14 (was data row now matching original row 10)
    $ws.Range("D14").Value = 44544
    $ws.Range("H14").Value = 'Inferno'
    $ws.Range("J14").Value = 12
    $ws.Range("K14").Value = 35000
    $ws.Range("L14").Value = 35000
    $ws.Range("M14").Value = 35000
    $ws.Range("P14").Value = 1400

    # Row 15 (was data row now matching original row 12)
    $ws.Range("D15").Value = 44460
    $ws.Range("J15").Value = 30
    $ws.Range("K15").Value = 95000
    $ws.Range("L15").Value = 95000
    $ws.Range("M15").Value = 95000
    $ws.Range("P15").Value = 3800

    # Row 16 (was data row now matching original row 4)
    $ws.Range("D16").Value = 44421
    $ws.Range("K16").Value = 75000
    $ws.Range("L16").Value = 75000
    $ws.Range("M16").Value = 75000
    $ws.Range("P16").Value = 3000

    # Row 17 (was data row now matching original row 16)
    $ws.Range("D17").Value = 44340
    $ws.Range("K17").Value = 35000
    $ws.Range("L17").Value = 35000
    $ws.Range("M17").Value = 35000
    $ws.Range("P17").Value = 1400

    # Row 18 (was data row now matching original row 14)
    $ws.Range("D18").Value = 44425
    $ws.Range("J18").Value = 15
    $ws.Range("K18").Value = 75000
    $ws.Range("L18").Value = 75000
    $ws.Range("M18").Value = 75000
    $ws.Range("P18").Value = 3000

Write-Output "Updated rows 2-18 with refreshed weekly price data."
